# "Corrected output location of reports from tests"
#
# The three criterion sheets (tabs 3, 4, 5 -- originally named
# "Criterion 2, Air Speed 0.1", "Criterion 1, Air Speed 0.1" and
# "Criterion 3, Air Speed 0.1") had been written out under the wrong
# labels/locations. This script rotates the underlying IES/MF result
# data among the three sheets, renames the sheets/table headers so that
# each sheet's label matches the data it now holds, and refreshes the
# small "readme" index table (sheet 1) to match: the Author/JobNo
# columns had been swapped, and the report date is bumped by a day.

$wb = $excel.ActiveWorkbook

$sheetCriterionA = $wb.Worksheets.Item(3)   # was "Criterion 2, Air Speed 0.1"
$sheetCriterionB = $wb.Worksheets.Item(4)   # was "Criterion 1, Air Speed 0.1"
$sheetCriterionC = $wb.Worksheets.Item(5)   # was "Criterion 3, Air Speed 0.1"

# ---------------------------------------------------------------------
# 1. Snapshot the C:D (IES Results / MF Results) data + the row 21/24
#    "F" column quirk for each of the three sheets BEFORE any writes,
#    since the sheets are about to trade data with each other.
# ---------------------------------------------------------------------
function Read-CriterionData($sheet) {
    $rows = @{}
    foreach ($r in 2..28) {
        $cText = $sheet.Cells.Item($r, 3).Text
        $dText = $sheet.Cells.Item($r, 4).Text
        $fText = $sheet.Cells.Item($r, 6).Text
        $rows[$r] = @{ C = $cText; D = $dText; F = $fText }
    }
    return $rows
}

$dataA = Read-CriterionData $sheetCriterionA
$dataB = Read-CriterionData $sheetCriterionB
$dataC = Read-CriterionData $sheetCriterionC

function Write-CriterionData($sheet, $rows) {
    foreach ($r in 2..28) {
        $row = $rows[$r]
        $sheet.Cells.Item($r, 3).Value = [double]$row.C
        $sheet.Cells.Item($r, 4).Value = [double]$row.D
        if ($row.F -eq "") {
            $sheet.Cells.Item($r, 6).ClearContents()
        } else {
            $sheet.Cells.Item($r, 6).Value = [double]$row.F
        }
    }
}

# Rotation: A <- B, B <- C, C <- A
Write-CriterionData $sheetCriterionA $dataB
Write-CriterionData $sheetCriterionB $dataC
Write-CriterionData $sheetCriterionC $dataA

# ---------------------------------------------------------------------
# 2. Rename the sheet tabs so each carries the label matching the data
#    it now holds. Do the 3-cycle through a temporary name to dodge
#    "sheet already exists" collisions.
# ---------------------------------------------------------------------
$sheetCriterionA.Name = "TEMP_SHEET_RENAME"
$sheetCriterionC.Name = "Criterion 2, Air Speed 0.1"
$sheetCriterionB.Name = "Criterion 3, Air Speed 0.1"
$sheetCriterionA.Name = "Criterion 1, Air Speed 0.1"

# ---------------------------------------------------------------------
# 3. Refresh the "Criterion N Absolute/Relative Change (%)" column
#    headers (also drives the Table column names) on each sheet so
#    they match the sheet's new identity.
# ---------------------------------------------------------------------
$sheetCriterionA.Range("E1").Value = "Criterion 1 Absolute Change"
$sheetCriterionA.Range("F1").Value = "Criterion 1 Relative Change (%)"

$sheetCriterionB.Range("E1").Value = "Criterion 3 Absolute Change"
$sheetCriterionB.Range("F1").Value = "Criterion 3 Relative Change (%)"

$sheetCriterionC.Range("E1").Value = "Criterion 2 Absolute Change"
$sheetCriterionC.Range("F1").Value = "Criterion 2 Relative Change (%)"

# ---------------------------------------------------------------------
# 4. readme sheet: the "JobNo" / "Author" columns were swapped, and the
#    report date moved from 20220614 to 20220615.
# ---------------------------------------------------------------------
$readme = $wb.Worksheets.Item(1)

$readme.Range("B1").Value = "Author"
$readme.Range("E1").Value = "JobNo"

foreach ($r in 2..5) {
    $readme.Cells.Item($r, 2).Value = "jovyan"   # Author
    $readme.Cells.Item($r, 5).Value = "/c/e"     # JobNo
    # Keep the date a literal text value (not auto-converted to a number).
    $readme.Cells.Item($r, 4).Value = "'20220615"
}

# sheet_name column (C) lists the criteria tabs in their new order.
$readme.Range("C3").Value = "Criterion 1, Air Speed 0.1"
$readme.Range("C4").Value = "Criterion 3, Air Speed 0.1"
$readme.Range("C5").Value = "Criterion 2, Air Speed 0.1"
